$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Hora) changes uniformly from 21 to 22 for every data row (2-51).
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "22"

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "293.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.24%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.61%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.949"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.31%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07328"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.03%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.275"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "26.22%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.730"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.74%"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.740"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.13%"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9085"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.06%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1684"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.80%"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08070"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "7.85%"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08150"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.70%"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03109"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.98%"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1008"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.80%"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001519"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.77%"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005727"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.11%"

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.31%"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.077"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.19%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1292"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.47%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.983"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-8.92%"

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.17%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04557"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.74%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001212"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.11%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004337"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.80%"

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.05%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01598"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.21%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04440"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.42%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007360"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.02%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.008666"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1325"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.33%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001937"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.20%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009522"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.24%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005970"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.31%"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.14%"

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "2.14%"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.14%"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.14%"
